$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: 'Stuck in the Moment' / 'Horn Glue'
$ws.Range("H40").Value = 1785.5385
$ws.Range("J40").Value = 1502.6
$ws.Range("L40").Value = 1502.6
$ws.Range("N40").Value = -1852.6

# Row 121: 'Mindful Medicine' / 'Tincture of Mind'
$ws.Range("H121").Value = 1065
$ws.Range("I121").Value = 797.5
$ws.Range("J121").Value = 1600
$ws.Range("K121").Value = 2392.5
$ws.Range("L121").Value = 4800
$ws.Range("M121").Value = -645.5
$ws.Range("N121").Value = -8294

# Row 132: 'Fast-forwarding Flora' / 'Growth Formula Lambda'
$ws.Range("H132").Value = 5162.364
$ws.Range("I132").Value = 6422.5
$ws.Range("J132").Value = 1802
$ws.Range("K132").Value = 19267.5
$ws.Range("L132").Value = 5406
$ws.Range("M132").Value = -16737.5
$ws.Range("N132").Value = -10466

# Row 138: 'All-night Crafting' / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 3779.513
$ws.Range("I138").Value = 3136.75
$ws.Range("J138").Value = 3945.3872
$ws.Range("K138").Value = 9410.25
$ws.Range("L138").Value = 11836.1616
$ws.Range("M138").Value = -4270.25
$ws.Range("N138").Value = -22116.1616

$ws = $wb.Worksheets.Item("ARM")
# Row 74: 'As the Bolt Flies' / 'Titanium Nugget'
$ws.Range("H74").Value = 1630.3529
$ws.Range("I74").Value = 1866.909
$ws.Range("J74").Value = 1196.6666
$ws.Range("K74").Value = 1866.909
$ws.Range("L74").Value = 1196.6666
$ws.Range("M74").Value = -992.9090000000001
$ws.Range("N74").Value = -2944.6666

# Row 77: 'Heavy Metal Banned (L)' / 'Titanium Nugget'
$ws.Range("H77").Value = 1630.3529
$ws.Range("I77").Value = 1866.909
$ws.Range("J77").Value = 1196.6666
$ws.Range("K77").Value = 9334.545
$ws.Range("L77").Value = 5983.333000000001
$ws.Range("M77").Value = -4966.545
$ws.Range("N77").Value = -14719.333

# Row 132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot'
$ws.Range("H132").Value = 2985.4
$ws.Range("I132").Value = 1495.9375
$ws.Range("J132").Value = 4239.684
$ws.Range("K132").Value = 4487.8125
$ws.Range("L132").Value = 12719.052
$ws.Range("M132").Value = -1957.8125
$ws.Range("N132").Value = -17779.052

$ws = $wb.Worksheets.Item("BSM")
# Row 134: 'Ruthenium Supremium' / 'Ruthenium Ingot'
$ws.Range("H134").Value = 2739.8462
$ws.Range("I134").Value = 1319.25
$ws.Range("K134").Value = 3957.75
$ws.Range("M134").Value = -1422.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 2281.3513
$ws.Range("I31").Value = 1968
$ws.Range("J31").Value = 3022
$ws.Range("K31").Value = 1968
$ws.Range("L31").Value = 3022
$ws.Range("M31").Value = -1673
$ws.Range("N31").Value = -3612

# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 2281.3513
$ws.Range("I34").Value = 1968
$ws.Range("J34").Value = 3022
$ws.Range("K34").Value = 1968
$ws.Range("L34").Value = 3022
$ws.Range("M34").Value = -1766
$ws.Range("N34").Value = -3426

# Row 99: 'O Pine' / 'Pine Lumber'
$ws.Range("H99").Value = 9002.923000000001
$ws.Range("I99").Value = 1470.6666
$ws.Range("J99").Value = 15459.143
$ws.Range("K99").Value = 1470.6666
$ws.Range("L99").Value = 15459.143
$ws.Range("M99").Value = 27.33339999999998
$ws.Range("N99").Value = -18455.143

# Row 107: 'Built to Last' / 'White Oak Lumber'
$ws.Range("H107").Value = 498.5
$ws.Range("I107").Value = 319.76923
$ws.Range("J107").Value = 830.4286
$ws.Range("K107").Value = 319.76923
$ws.Range("L107").Value = 830.4286
$ws.Range("M107").Value = 1600.23077
$ws.Range("N107").Value = -4670.4286

# Row 126: 'A Better Conductor' / 'Red Pine Lumber'
$ws.Range("H126").Value = 9002.923000000001
$ws.Range("I126").Value = 1470.6666
$ws.Range("J126").Value = 15459.143
$ws.Range("K126").Value = 4411.9998
$ws.Range("L126").Value = 46377.429
$ws.Range("M126").Value = -1941.9998
$ws.Range("N126").Value = -51317.429

# Row 132: 'Hull Lotta Damage' / 'Ginseng Lumber'
$ws.Range("H132").Value = 2597.238
$ws.Range("I132").Value = 2577.182
$ws.Range("J132").Value = 2619.3
$ws.Range("K132").Value = 7731.545999999999
$ws.Range("L132").Value = 7857.900000000001
$ws.Range("M132").Value = -5201.545999999999
$ws.Range("N132").Value = -12917.9

# Row 134: 'Wood You Be Quiet' / 'Ceiba Lumber'
$ws.Range("H134").Value = 4666.8
$ws.Range("I134").Value = 5509
$ws.Range("J134").Value = 3403.5
$ws.Range("K134").Value = 16527
$ws.Range("L134").Value = 10210.5
$ws.Range("M134").Value = -13992
$ws.Range("N134").Value = -15280.5

$ws = $wb.Worksheets.Item("CUL")
# Row 88: "Don't Let It Fall Apart" / 'Liver-cheese Sandwich'
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91: 'Better Come Back with a Sandwich (L)' / 'Liver-cheese Sandwich'
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("N91").ClearContents()

# Row 131: 'The Mountain Steeped' / 'Tsai tou Vounou'
$ws.Range("H131").Value = 2951.25
$ws.Range("J131").Value = 3154.423
$ws.Range("L131").Value = 9463.269
$ws.Range("N131").Value = -19543.269

# Row 132: 'More Mezcal' / 'Cooking Mezcal'
$ws.Range("H132").Value = 1430.4
$ws.Range("I132").Value = 1411.8235
$ws.Range("J132").Value = 1469.875
$ws.Range("K132").Value = 12706.4115
$ws.Range("L132").Value = 13228.875
$ws.Range("M132").Value = -10176.4115
$ws.Range("N132").Value = -18288.875

$ws = $wb.Worksheets.Item("GSM")
# Row 107: 'Whetstones for the Workers' / 'Hard Mudstone Whetstone'
$ws.Range("H107").Value = 1470.6666
$ws.Range("I107").Value = 1801.8334
$ws.Range("J107").Value = 146
$ws.Range("K107").Value = 1801.8334
$ws.Range("L107").Value = 146
$ws.Range("M107").Value = 118.1666
$ws.Range("N107").Value = -3986

# Row 132: 'On Board for Lar' / 'Lar Ingot'
$ws.Range("H132").Value = 4864.25
$ws.Range("I132").Value = 3224.1667
$ws.Range("K132").Value = 9672.500100000001
$ws.Range("M132").Value = -7142.500100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 38: 'Emergency Patches' / 'Skull Eyepatch'
$ws.Range("H38").Value = 22800
$ws.Range("J38").Value = 22800
$ws.Range("L38").Value = 22800
$ws.Range("N38").Value = -23620

# Row 108: 'Girding for Glory' / 'Smilodonskin Trousers of Maiming'
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").ClearContents()

# Row 136: "Respect for Br'aax" / "Br'aax Leather"
$ws.Range("H136").Value = 25646504
$ws.Range("I136").Value = 7200.8
$ws.Range("J136").Value = 41671068
$ws.Range("K136").Value = 21602.4
$ws.Range("L136").Value = 125013204
$ws.Range("M136").Value = -19052.4
$ws.Range("N136").Value = -125018304

$ws = $wb.Worksheets.Item("WVR")
# Row 64: 'Ribbon of Remembrance' / 'Rainbow Ribbon of Healing'
$ws.Range("H64").Value = 29914
$ws.Range("J64").Value = 29914
$ws.Range("L64").Value = 29914
$ws.Range("N64").Value = -30410

# Row 67: 'The Road Was a Ribbon of Moonlight (L)' / 'Rainbow Ribbon of Healing'
$ws.Range("H67").Value = 29914
$ws.Range("J67").Value = 29914
$ws.Range("L67").Value = 29914
$ws.Range("N67").Value = -31630

# Row 132: 'Comfy Cabins' / 'Snow Cotton Cloth'
$ws.Range("H132").Value = 3700.5386
$ws.Range("I132").Value = 3467
$ws.Range("J132").Value = 4226
$ws.Range("K132").Value = 10401
$ws.Range("L132").Value = 12678
$ws.Range("M132").Value = -7871
$ws.Range("N132").Value = -17738
